# Daily attendance processing - 2025-10-29 22:46:35
# Update "Recorded By" (column G) values to reflect reordered recorder names
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value = "backup@backdoor.com, System, system"
$ws.Range("G4").Value = "backup@backdoor.com, System"
$ws.Range("G5").Value = "backup@backdoor.com, System"
$ws.Range("G8").Value = "backup@backdoor.com, System"
$ws.Range("G11").Value = "System, dnasr281@gmail.com"
$ws.Range("G17").Value = "System, dnasr281@gmail.com"
$ws.Range("G29").Value = "backup@backdoor.com, System, system"
$ws.Range("G31").Value = "backup@backdoor.com, System"
$ws.Range("G32").Value = "backup@backdoor.com, System"
$ws.Range("G35").Value = "backup@backdoor.com, System"
$ws.Range("G38").Value = "System, dnasr281@gmail.com"
$ws.Range("G44").Value = "System, dnasr281@gmail.com"
$ws.Range("G56").Value = "backup@backdoor.com, System, system"
$ws.Range("G58").Value = "backup@backdoor.com, System"
$ws.Range("G59").Value = "backup@backdoor.com, System"
$ws.Range("G62").Value = "backup@backdoor.com, System"
$ws.Range("G65").Value = "System, dnasr281@gmail.com"
$ws.Range("G71").Value = "System, dnasr281@gmail.com"
$ws.Range("G83").Value = "backup@backdoor.com, System"
$ws.Range("G84").Value = "backup@backdoor.com, System"
$ws.Range("G85").Value = "backup@backdoor.com, System"
$ws.Range("G96").Value = "System, dnasr281@gmail.com"
$ws.Range("G97").Value = "System, dnasr281@gmail.com"
$ws.Range("G99").Value = "System, dnasr281@gmail.com"
$ws.Range("G109").Value = "backup@backdoor.com, System"
$ws.Range("G110").Value = "backup@backdoor.com, System"
$ws.Range("G111").Value = "backup@backdoor.com, System"
$ws.Range("G122").Value = "System, dnasr281@gmail.com"
$ws.Range("G123").Value = "System, dnasr281@gmail.com"
$ws.Range("G125").Value = "System, dnasr281@gmail.com"
$ws.Range("G135").Value = "backup@backdoor.com, System"
$ws.Range("G136").Value = "backup@backdoor.com, System"
$ws.Range("G137").Value = "backup@backdoor.com, System"
$ws.Range("G148").Value = "System, dnasr281@gmail.com"
$ws.Range("G149").Value = "System, dnasr281@gmail.com"
$ws.Range("G151").Value = "System, dnasr281@gmail.com"
